$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows: small recalculated statistics (H column, and G
# column for rows 4 & 7) produced by re-running the comparison with the
# extra "B" variable folded in ---
$ws.Range("H2").Value  = 0.9646803381026858

$ws.Range("H3").Value  = 0.003427539896620487

$ws.Range("G4").Value  = 0.09289560375421288
$ws.Range("H4").Value  = 0.9261800894161204

$ws.Range("H5").Value  = 0.7308425919715824

$ws.Range("H6").Value  = 0.6529291562813258

$ws.Range("G7").Value  = -3.810620297534498
$ws.Range("H7").Value  = 0.0002446487635881599

$ws.Range("H8").Value  = 0.0001743429192672115

$ws.Range("H9").Value  = 0.2505463453197566

$ws.Range("H10").Value = 0.4877664561713912

# --- Append the new "B" variable row ---
$ws.Range("A11").Value = "B"
$ws.Range("B11").Value = -0.5685414527331637
$ws.Range("C11").Value = -2.810999543124124
$ws.Range("D11").Value = 0.7609251830041186
$ws.Range("E11").Value = 17.17363614673463
$ws.Range("F11").Value = 0.2964689350684678
$ws.Range("G11").Value = 0.9131334187456903
$ws.Range("H11").Value = 0.3634599521358221
$ws.Range("I11").Value = "No"
